$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Update the cached date field text on every slide (date placeholders)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
}

$conn = $s.Shapes.Item("Straight Arrow Connector 2")
$conn.Line.EndArrowheadStyle = 2
